# Apply updated crypto price/volume data to Sheet1 (columns D and E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "60.456.21"
$ws.Cells.Item(2, 5).Value = "  +0.43%  "

$ws.Cells.Item(3, 4).Value = "2.599.11"
$ws.Cells.Item(3, 5).Value = "  +0.78%  "

$ws.Cells.Item(4, 5).Value = "  -0.18%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "514.23"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = "  +2.03%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "153.30"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = "  -0.70%  "

$ws.Cells.Item(7, 5).Value = "  +0.12%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.599"
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = "  +3.49%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "6.63"
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = "  +0.76%  "

$ws.Cells.Item(10, 5).Value = "  +0.35%  "

$ws.Cells.Item(11, 5).Value = "  +0.93%  "

$ws.Cells.Item(12, 5).Value = "  +1.80%  "

$ws.Cells.Item(13, 4).Value = "3.054.52"
$ws.Cells.Item(13, 5).Value = "  +0.61%  "

$ws.Cells.Item(14, 4).Value = "60.495.50"
$ws.Cells.Item(14, 5).Value = "  +0.41%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "21.57"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(15, 5).Value = "  -0.27%  "

$ws.Cells.Item(16, 5).Value = "  +0.61%  "

$ws.Cells.Item(17, 4).Value = "2.605.22"
$ws.Cells.Item(17, 5).Value = "  -0.10%  "

$ws.Cells.Item(18, 5).Value = "  -1.05%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "358.42"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = "  +3.80%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "10.54"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = "  +2.00%  "

$ws.Cells.Item(21, 5).Value = "  +1.99%  "

$ws.Cells.Item(22, 5).Value = "  +0.15%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "61.07"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = "  +1.38%  "

$ws.Cells.Item(24, 5).Value = "  +1.30%  "

$ws.Cells.Item(25, 4).Value = "2.716.88"
$ws.Cells.Item(25, 5).Value = "  +0.39%  "

$ws.Cells.Item(26, 5).Value = "  +0.34%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0.999"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = "  +0.66%  "

$ws.Cells.Item(28, 4).Value = "0.0₃0833"
$ws.Cells.Item(28, 5).Value = "  -2.22%  "

$ws.Cells.Item(29, 5).Value = "  -2.71%  "

$ws.Cells.Item(30, 5).Value = "  +0.16%  "

$ws.Cells.Item(31, 5).Value = "  +0.71%  "

$ws.Cells.Item(32, 5).Value = "  +1.95%  "

$ws.Cells.Item(33, 5).Value = "  +3.74%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "150.21"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = "  -3.71%  "

$ws.Cells.Item(35, 5).Value = "  +0.43%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.918"
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = "  +7.67%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.18"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = "  -1.11%  "

$ws.Cells.Item(38, 5).Value = "  +0.54%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "36.29"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = "  +1.57%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.841"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = "  -0.49%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "3.74"
$ws.Cells.Item(41, 4).ClearFormats()

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "288.18"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = "  -2.65%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.102"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = "  +1.84%  "

$ws.Cells.Item(44, 5).Value = "  -0.62%  "

$ws.Cells.Item(45, 5).Value = "  -0.04%  "

$ws.Cells.Item(46, 5).Value = "  -2.35%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "19.56"
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = "  -1.06%  "

$ws.Cells.Item(48, 5).Value = "  -0.27%  "

$ws.Cells.Item(49, 5).Value = "  +0.50%  "

$ws.Cells.Item(50, 5).Value = "  +0.36%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "19.20"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = "  +9.78%  "
